$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Warning Signs" song entry (row 2), shifting subsequent rows up.
$ws.Rows.Item(2).Delete()

$ws.Range("A8").Select()
